$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.232.70"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "1.787.57"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.96"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("E6").Value = "  +1.34%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("E10").Value = "  -0.36%  "
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").Value = "2.045.55"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").Value = "1.789.11"
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.00"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.67%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.625"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "34.204.70"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.18"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.00"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.99"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").Value = "0.0₃0798"
$ws.Range("E20").Value = "  +2.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.93"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("E23").Value = "  +0.88%  "
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.53"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.17"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.32"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("E30").Value = "  -0.88%  "
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("E32").Value = "  +2.13%  "
$ws.Range("E33").Value = "  +3.44%  "
$ws.Range("E34").Value = "  -1.59%  "
$ws.Range("D35").Value = "1.438.24"
$ws.Range("E35").Value = "  +0.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.59"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +9.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.663"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.10%  "
$ws.Range("E38").Value = "  +1.10%  "
$ws.Range("E39").Value = "  -1.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "81.72"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.36%  "
$ws.Range("E41").Value = "  +5.45%  "
$ws.Range("E42").Value = "  +1.43%  "
$ws.Range("E43").Value = "  +1.51%  "
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0520"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.13%  "
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "1.942.00"
$ws.Range("E48").Value = "  -0.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.36"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.11%  "
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("E51").Value = "  -5.90%  "
